$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sector names (column A) and average correlation values (column B) for rows 2-48
$data = @(
    ,@("Construction Materials(8)", 0.6008465032473096)
    ,@("Gas Utilities(12)", 0.5186801690767942)
    ,@("Multi-Utilities(18)", 0.500909479822461)
    ,@("Marine(15)", 0.4852696261099415)
    ,@("Construction & Engineering(21)", 0.4515538674136891)
    ,@("Mortgage Real Estate Investment Trust...(16)", 0.4474746057728351)
    ,@("Electric Utilities(28)", 0.4335275976640242)
    ,@("Wireless Telecommunication Services(14)", 0.3941202375365044)
    ,@("Containers & Packaging(12)", 0.380250765115083)
    ,@("Equity Real Estate Investment Trusts ...(98)", 0.3782152918493715)
    ,@("Building Products(24)", 0.3715054006124927)
    ,@("ETF(303)", 0.3603432939999822)
    ,@("Auto Components(21)", 0.3381824331755282)
    ,@("Energy Equipment & Services(38)", 0.3343197033573453)
    ,@("Capital Markets(76)", 0.328547175184346)
    ,@("Life Sciences Tools & Services(19)", 0.3191484823244686)
    ,@("Chemicals(52)", 0.3159253043852266)
    ,@("Oil, Gas & Consumable Fuels(125)", 0.300488469143191)
    ,@("Machinery(86)", 0.2914865978817048)
    ,@("Insurance(75)", 0.2796246493324812)
    ,@("Diversified Consumer Services(17)", 0.2779035428496974)
    ,@("Diversified Telecommunication Services(20)", 0.2711477173287575)
    ,@("Pharmaceuticals(53)", 0.2653524010583758)
    ,@("Professional Services(35)", 0.2566271024674783)
    ,@("Technology Hardware, Storage & Periph...(19)", 0.2554864017417568)
    ,@("Banks(251)", 0.2475183934629234)
    ,@("Metals & Mining(106)", 0.2401803739102153)
    ,@("Communications Equipment(45)", 0.2349719760356817)
    ,@("Real Estate Management & Development(23)", 0.2314824023157566)
    ,@("Media(42)", 0.2298310009877373)
    ,@("Household Durables(39)", 0.2268738898233296)
    ,@("Trading Companies & Distributors(25)", 0.2261234752372032)
    ,@("IT Services(52)", 0.2229508408562265)
    ,@("Health Care Providers & Services(47)", 0.2137223389562678)
    ,@("Semiconductors & Semiconductor Equipment(70)", 0.213702293559098)
    ,@("Health Care Equipment & Supplies(86)", 0.2049337064446754)
    ,@("Road & Rail(22)", 0.199812233398616)
    ,@("Thrifts & Mortgage Finance(47)", 0.184373386344324)
    ,@("Textiles, Apparel & Luxury Goods(29)", 0.178778919123491)
    ,@("Software(70)", 0.17000043990075)
    ,@("Commercial Services & Supplies(52)", 0.1688693416197126)
    ,@("Food Products(46)", 0.1671951914030452)
    ,@("Hotels, Restaurants & Leisure(51)", 0.1528489294252231)
    ,@("Aerospace & Defense(37)", 0.1492547521842049)
    ,@("Specialty Retail(59)", 0.1386370704968179)
    ,@("Electronic Equipment, Instruments & C...(78)", 0.1331343574556911)
    ,@("Biotechnology(128)", 0.09617365007536748)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = [double]$data[$i][1]
}

